# Updated grades for 4, 5, 6 UE assignments; added a comment about the
# missing "shooting arrow" action for assignment 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Block "4 UE" (columns R:T) ---
$ws.Range("R2").Value = "4 UE"
$ws.Range("S2").Value = 100
$ws.Range("T2").Value = "Good Job!"

# --- Block "5 UE" (columns V:X) ---
$ws.Range("V2").Value = "5 UE"
$ws.Range("W2").Value = 50
$ws.Range("X2").Value = "You didn’t implement a new action for shooting arrow 2"

# --- Block "6 UE" (columns Z:AB) ---
$ws.Range("Z2").Value = "6 UE"
$ws.Range("AA2").Value = 100
$ws.Range("AB2").Value = "Good Job!"

# The comment text in X2 is long, so Excel grows row 2 to show it in full.
$ws.Rows.Item(2).RowHeight = 90

# Column X (24) is widened to comfortably fit the new, longer comment text.
$ws.Columns.Item(24).ColumnWidth = 19.42578125

# Reflect the newly-entered grade block as the current selection, matching
# what a user would see right after typing the values in.
$ws.Range("R2:AB2").Select()
